$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 29.31851203064848
$ws.Range("B3").Value = 26.3675529734944
$ws.Range("B4").Value = 20.70803954796096
$ws.Range("B5").Value = 10.50027259277895
$ws.Range("B6").Value = 7.862753710051312
$ws.Range("B7").Value = 5.242869145065904
